$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.781.37'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.746.08'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.69'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5081'
$ws.Range('E7').Value = '  +3.19%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2670'
$ws.Range('E8').Value = '  +6.13%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06177'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').Value = '1.742.15'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06947'
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.26'
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6210'
$ws.Range('E13').Value = '  +9.84%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.467'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '77.68'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '25.801.75'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.58'
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000006629'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '1.967.66'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.042'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.247'
$ws.Range('E23').Value = '  +4.79%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.123'
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '136.50'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.462'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.97'
$ws.Range('E27').Value = '  +2.28%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.765'
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '102.41'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.08128'
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.686'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.384'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04395'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.647'
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9926'
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6015'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.593'
$ws.Range('E37').Value = '  -3.06%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01554'
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.933'
$ws.Range('E39').Value = '  -3.70%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('B41').Value = 'PaxosStandard'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.55'
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3816'
$ws.Range('E43').Value = '  +2.53%  '
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.898'
$ws.Range('E45').Value = '  -5.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.05497'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1094'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.909'
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '30.01'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '52.41'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.004'
$ws.Range('E51').Value = '  +0.45%  '
